# Update deployment plot labels
# - Update B11:B18 with new cumulative capacity values
# - Remove the now-obsolete last row (previously row 19, year 2048)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cumulative capacity values for years 2040-2046 (rows 11-17)
$ws.Range("B11").Value = 26541.71939820418
$ws.Range("B12").Value = 31180.7898679445
$ws.Range("B13").Value = 35927.35275193814
$ws.Range("B14").Value = 40590.58779680728
$ws.Range("B15").Value = 45096.14405351903
$ws.Range("B16").Value = 49796.00380965499
$ws.Range("B17").Value = 53747.17892691874

# Row 18 (year 2047) now takes on the value previously held by row 19 (year 2048)
$ws.Range("B18").Value = 54601.2927938712

# Delete the now-redundant last row (year 2048), shifting rows up
$ws.Rows("19").Delete()
